# Changes of 28th july 2022
#
# The PackageTrackNum / ShipmentTrackNum columns (C and D) of Sheet1 get a
# fresh batch of FedEx tracking numbers. Column C (rows 2-22) always holds a
# value, and column D mirrors C on the rows where a "ShipmentTrackNum" is
# also recorded (rows 5-7 and 13-17).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These 12-digit values look numeric, so Excel would otherwise store them as
# plain numbers. Force the cells to stay text (as they originally were,
# stored as shared strings) before writing the new values.
$ws.Range("C2:C22").NumberFormat = "@"
$ws.Range("D5:D7").NumberFormat = "@"
$ws.Range("D13:D17").NumberFormat = "@"

$ws.Range("C2").Value = "320018767744"
$ws.Range("C3").Value = "320018767814"
$ws.Range("C4").Value = "320018767847"
$ws.Range("C5").Value = "320018767869"
$ws.Range("C6").Value = "320018767906"
$ws.Range("C7").Value = "320018767928"
$ws.Range("C8").Value = "320018768155"
$ws.Range("C9").Value = "320018768177"
$ws.Range("C10").Value = "320018768203"
$ws.Range("C11").Value = "320018768225"
$ws.Range("C12").Value = "320018768269"
$ws.Range("C13").Value = "320018768280"
$ws.Range("C14").Value = "320018768317"
$ws.Range("C15").Value = "320018768339"
$ws.Range("C16").Value = "320018768361"
$ws.Range("C17").Value = "320018768383"
$ws.Range("C18").Value = "320018768420"
$ws.Range("C19").Value = "320018768442"
$ws.Range("C20").Value = "320018768475"
$ws.Range("C21").Value = "320018768497"
$ws.Range("C22").Value = "320018768523"

$ws.Range("D5").Value = "320018767869"
$ws.Range("D6").Value = "320018767906"
$ws.Range("D7").Value = "320018767928"
$ws.Range("D13").Value = "320018768280"
$ws.Range("D14").Value = "320018768317"
$ws.Range("D15").Value = "320018768339"
$ws.Range("D16").Value = "320018768361"
$ws.Range("D17").Value = "320018768383"

# Put the number format back to General so the cells end up exactly like
# the originals: plain text (shared-string) cells with no custom style.
$ws.Range("C2:C22").Style = "Normal"
$ws.Range("D5:D7").Style = "Normal"
$ws.Range("D13:D17").Style = "Normal"
